$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename custom field header columns (S1: CF1 -> Custom Field 1, T1: CF2 -> Custom field   2)
# Set T1 first so the new shared strings land in the same append order as the target file.
$ws.Range("T1").Value = "Custom field   2"
$ws.Range("S1").Value = "Custom Field 1"

# Update the active selection to S2 (was T6)
$ws.Range("S2").Select() | Out-Null
